$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a value to be stored as literal TEXT (no Excel smart-typing
# into a number/date/percentage), while leaving the cell's number format
# untouched (the leading apostrophe otherwise tags the cell with a
# "quote prefix" style, so we strip that again with ClearFormats).
function Set-TextValue($addr, $text) {
    $ws.Range($addr).Value = "'" + $text
    $ws.Range($addr).ClearFormats()
}

# Helper: force a value to be stored as an empty literal TEXT cell (present
# in the sheet, but with no characters) instead of being removed outright
# (plain `.Value = ""` deletes the cell entirely).
function Set-EmptyTextValue($addr) {
    $ws.Range($addr).Value = "'"
    $ws.Range($addr).ClearFormats()
}

function Set-NumberValue($addr, $number) {
    $ws.Range($addr).Value = $number
}

function Set-DataRow($r, $a, $b, $c, $d, $e, $f, $g, $h, $i, $j, $k, $cAsText) {
    Set-TextValue "A$r" $a
    Set-TextValue "B$r" $b
    if ($cAsText) {
        Set-TextValue "C$r" $c
    } else {
        Set-NumberValue "C$r" $c
    }
    Set-TextValue "D$r" $d
    Set-TextValue "E$r" $e
    Set-TextValue "F$r" $f
    Set-TextValue "G$r" $g
    Set-TextValue "H$r" $h
    if ($i -eq "") {
        Set-EmptyTextValue "I$r"
    } else {
        Set-TextValue "I$r" $i
    }
    Set-TextValue "J$r" $j
    Set-NumberValue "K$r" $k
}

# --- Existing rows 25-33: project_id (C) becomes numeric, test (I) N/A -> blank ---
Set-NumberValue "C25" 139863102765096
Set-NumberValue "C26" 139863102765096
Set-NumberValue "C27" 139863102765096
Set-NumberValue "C28" 139863258720298
Set-NumberValue "C29" 139863258720298
Set-NumberValue "C30" 139863258720298
Set-NumberValue "C31" 139863342313515
Set-NumberValue "C32" 139863342313515
Set-NumberValue "C33" 139863342313515

Set-EmptyTextValue "I25"
Set-EmptyTextValue "I26"
Set-EmptyTextValue "I27"
Set-EmptyTextValue "I28"
Set-EmptyTextValue "I29"
Set-EmptyTextValue "I30"
Set-EmptyTextValue "I31"
Set-EmptyTextValue "I32"
# I33 is unchanged (stays "N/A")

# --- New rows 34-42 (test https errors) ---
Set-DataRow 34 "Epochs" "object-detection" "139863102765096" "7.02 minutes" "RtmDet-[9M]" "{'height': 600, 'width': 600, 'paddingValue': 0}" "99%" "93%" "N/A" "2025-06-10 11:20:27" 50 $true
Set-DataRow 35 "Epochs" "object-detection" "139863102765096" "13.59 minutes" "RepPoints-[20M]" "{'height': 600, 'width': 600, 'paddingValue': 0}" "100%" "95%" "N/A" "2025-06-10 11:36:13" 100 $true
Set-DataRow 36 "Epochs" "object-detection" "139863102765096" "15.58 minutes" "RepPoints-[37M]" "{'height': 600, 'width': 600, 'paddingValue': 0}" "100%" "96%" "N/A" "2025-06-10 11:54:32" 150 $true
Set-DataRow 37 "Epochs" "classification" "139863258720298" "1.24 minutes" "ConvNext-[29M]" "{'height': 512, 'width': 512, 'paddingValue': 0}" "100%" "100%" "N/A" "2025-06-10 11:57:10" 50 $true
Set-DataRow 38 "Epochs" "classification" "139863258720298" "0.73 minutes" "ConvNext-[16M]" "{'height': 512, 'width': 512, 'paddingValue': 0}" "100%" "100%" "N/A" "2025-06-10 12:00:02" 100 $true
Set-DataRow 39 "Epochs" "classification" "139863258720298" "0.82 minutes" "ConvNext-[16M]" "{'height': 512, 'width': 512, 'paddingValue': 0}" "100%" "100%" "N/A" "2025-06-10 12:02:17" 150 $true
Set-DataRow 40 "Epochs" "segmentation" "139863342313515" "4.09 minutes" "SegFormer-[14M]" "{'height': 800, 'width': 800, 'paddingValue': 0}" "94%" "89%" "N/A" "2025-06-10 12:07:35" 50 $true
Set-DataRow 41 "Epochs" "segmentation" "139863342313515" "9.59 minutes" "FastVit-[14M]" "{'height': 800, 'width': 800, 'paddingValue': 0}" "96%" "94%" "N/A" "2025-06-10 12:18:05" 100 $true
Set-DataRow 42 "Epochs" "segmentation" "139863342313515" "12.06 minutes" "SegFormer-[14M]" "{'height': 800, 'width': 800, 'paddingValue': 0}" "95%" "91%" "N/A" "2025-06-10 12:31:35" 150 $true

Write-Host "edit applied"
